# Add a new data row (row 63) to the CRM accuracy log, recording a run with
# the CRM/pH buffers opened 11/8/2019 (matches the new shared string + the
# trailing row added in the target sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Date (column A) ---------------------------------------------------
# Set the serial date value, then copy just the number format from the
# previous date cell (A62) so this cell reuses the existing date style
# instead of Excel/iron_native registering a brand-new number format.
$ws.Range("A63").Value = 43781
$ws.Range("A62").Copy() | Out-Null
$ws.Range("A63").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# --- CRM value / Batch value (columns B, C) -----------------------------
$ws.Range("B63").Value = 2147.3009674444802
$ws.Range("C63").Value = 2207.0300000000002

# --- % off (column D) ----------------------------------------------------
# Continues the same "100*(B-C)/C" series used by the rows above it.
$ws.Range("D63").Formula = "=100*(B63-C63)/C63"

# --- Batch # (column E) ---------------------------------------------------
$ws.Range("E63").Value = 169

# --- Notes (column F) ------------------------------------------------------
$ws.Range("F63").Value = "crm opened 11/8/2019"

# Move the active selection to the next empty notes cell, matching what
# Excel leaves selected after data entry on the last row.
$ws.Range("F64").Select() | Out-Null
